# Update stats for 2026-01 (row 26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6485
$ws.Range("C26").Value = 1010
$ws.Range("D26").Value = 6043044
$ws.Range("E26").Value = 931.849498843485
$ws.Range("F26").Value = 9.599459185397997
$ws.Range("G26").Value = 7.218683651804669
$ws.Range("H26").Value = 25.84659488832888
